$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.534.77'
$ws.Range('E2').Value = '  +4.28%  '
$ws.Range('D3').Value = '3.505.91'
$ws.Range('E3').Value = '  +3.87%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '585.32'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.97'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.72%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.479'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.38%  '
$ws.Range('E9').Value = '  +0.67%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.127'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.54%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.399'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.33%  '
$ws.Range('D12').Value = '4.115.09'
$ws.Range('E12').Value = '  +4.17%  '
$ws.Range('E13').Value = '  +7.55%  '
$ws.Range('E14').Value = '  -0.72%  '
$ws.Range('D15').Value = '3.518.15'
$ws.Range('E15').Value = '  +4.06%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000175'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.51%  '
$ws.Range('D17').Value = '63.608.42'
$ws.Range('E17').Value = '  +4.25%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.28'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.39%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.37'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '9.48'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +7.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '396.42'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.69%  '
$ws.Range('B22').Value = 'Polygon'
$ws.Range('C22').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.566'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.25%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '75.86'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.26%  '
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000121'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +7.36%  '
$ws.Range('D26').Value = '3.652.91'
$ws.Range('E26').Value = '  +3.96%  '
$ws.Range('E27').Value = '  +1.23%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.90'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +10.62%  '
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.30'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.22%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.17'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.60%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.45'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +8.63%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.99'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.36'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +9.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.19'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.80%  '
$ws.Range('B37').Value = 'EnergySwap'
$ws.Range('C37').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '32.34'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +26.94%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '173.09'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.29%  '
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.59'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +10.05%  '
$ws.Range('D40').Value = '3.550.72'
$ws.Range('E40').Value = '  +4.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0780'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.78%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.803'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.71%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.76'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +8.62%  '
$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.54'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.68%  '
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '42.27'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.40%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.22'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +10.26%  '
$ws.Range('D47').Value = '2.604.31'
$ws.Range('E47').Value = '  +6.47%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.95'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +7.36%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.24'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +8.93%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.79'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.96%  '
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0272'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.72%  '
